$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1371.2693
$ws.Range("J98").Value = 1438.1666
$ws.Range("I98").Value = 1313.9286
$ws.Range("K98").Value = 1313.9286
$ws.Range("M98").Value = 184.0714
$ws.Range("L98").Value = 1438.1666
$ws.Range("N98").Value = -4434.1666
$ws.Range("K107").Value = 491.2857
$ws.Range("M107").Value = 1428.7143
$ws.Range("H107").Value = 512.2895
$ws.Range("I107").Value = 491.2857
$ws.Range("K122").Value = 3941.7858
$ws.Range("M122").Value = -1491.7858
$ws.Range("L122").Value = 4314.4998
$ws.Range("N122").Value = -9214.4998
$ws.Range("H122").Value = 1371.2693
$ws.Range("J122").Value = 1438.1666
$ws.Range("I122").Value = 1313.9286
$ws.Range("H137").Value = 9805821
$ws.Range("J137").Value = 2071.4285
$ws.Range("I137").Value = 13515348
$ws.Range("M137").Value = -40543494
$ws.Range("L137").Value = 6214.2855
$ws.Range("K137").Value = 40546044
$ws.Range("N137").Value = -11314.2855
$ws.Range("I138").Value = 1389.5333
$ws.Range("K138").Value = 4168.5999
$ws.Range("M138").Value = 971.4000999999998
$ws.Range("L138").Value = 6402.7272
$ws.Range("N138").Value = -16682.7272
$ws.Range("H138").Value = 1704.6025
$ws.Range("J138").Value = 2134.2424

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1224.8462
$ws.Range("J2").Value = 2544.75
$ws.Range("I2").Value = 638.2222
$ws.Range("M2").Value = -525.2222
$ws.Range("L2").Value = 2544.75
$ws.Range("K2").Value = 638.2222
$ws.Range("N2").Value = -2770.75
$ws.Range("K32").Value = 5038.815
$ws.Range("M32").Value = -4751.815
$ws.Range("L32").Value = 23840180
$ws.Range("N32").Value = -23840754
$ws.Range("H32").Value = 3517586
$ws.Range("J32").Value = 23840180
$ws.Range("I32").Value = 5038.815
$ws.Range("M45").Value = -1098.3529
$ws.Range("L45").Value = 4847.3
$ws.Range("K45").Value = 1475.3529
$ws.Range("H45").Value = 2724.2222
$ws.Range("N45").Value = -5601.3
$ws.Range("J45").Value = 4847.3
$ws.Range("I45").Value = 1475.3529
$ws.Range("M61").Value = -890.4572000000001
$ws.Range("L61").Value = 3780.5293
$ws.Range("K61").Value = 1102.4572
$ws.Range("N61").Value = -4204.5293
$ws.Range("H61").Value = 1977.9807
$ws.Range("J61").Value = 3780.5293
$ws.Range("I61").Value = 1102.4572
$ws.Range("K110").Value = 1463.35
$ws.Range("M110").Value = 581.6500000000001
$ws.Range("L110").Value = 5670.8335
$ws.Range("N110").Value = -9760.833500000001
$ws.Range("H110").Value = 2434.3076
$ws.Range("J110").Value = 5670.8335
$ws.Range("I110").Value = 1463.35
$ws.Range("M116").Value = 1655.7778
$ws.Range("L116").Value = 2544.75
$ws.Range("K116").Value = 638.2222
$ws.Range("H116").Value = 1224.8462
$ws.Range("N116").Value = -7132.75
$ws.Range("J116").Value = 2544.75
$ws.Range("I116").Value = 638.2222
$ws.Range("K122").Value = 4748.4
$ws.Range("M122").Value = -2298.4
$ws.Range("L122").Value = 11373.9
$ws.Range("N122").Value = -16273.9
$ws.Range("H122").Value = 3055.1333
$ws.Range("J122").Value = 3791.3
$ws.Range("I122").Value = 1582.8
$ws.Range("M132").Value = -1468.5116
$ws.Range("L132").Value = 15644700
$ws.Range("K132").Value = 3998.5116
$ws.Range("N132").Value = -15649760
$ws.Range("H132").Value = 1868879.2
$ws.Range("J132").Value = 5214900
$ws.Range("I132").Value = 1332.8372
$ws.Range("K136").Value = 3307.3716
$ws.Range("M136").Value = -757.3716000000004
$ws.Range("L136").Value = 11341.5879
$ws.Range("N136").Value = -16441.5879
$ws.Range("H136").Value = 1977.9807
$ws.Range("J136").Value = 3780.5293
$ws.Range("I136").Value = 1102.4572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I3").Value = 638.2222
$ws.Range("K3").Value = 638.2222
$ws.Range("M3").Value = -524.2222
$ws.Range("L3").Value = 2544.75
$ws.Range("N3").Value = -2772.75
$ws.Range("H3").Value = 1224.8462
$ws.Range("J3").Value = 2544.75
$ws.Range("I86").Value = 1600.4546
$ws.Range("K86").Value = 1600.4546
$ws.Range("M86").Value = -477.4546
$ws.Range("L86").Value = 3007
$ws.Range("N86").Value = -5253
$ws.Range("H86").Value = 1717.6666
$ws.Range("J86").Value = 3007
$ws.Range("L87").Value = 46000
$ws.Range("N87").Value = -48496
$ws.Range("H87").Value = 46000
$ws.Range("J87").Value = 46000
$ws.Range("I89").Value = 1600.4546
$ws.Range("M89").Value = -2386.273
$ws.Range("K89").Value = 8002.273
$ws.Range("L89").Value = 15035
$ws.Range("N89").Value = -26267
$ws.Range("H89").Value = 1717.6666
$ws.Range("J89").Value = 3007
$ws.Range("L90").Value = 138000
$ws.Range("H90").Value = 46000
$ws.Range("N90").Value = -150480
$ws.Range("J90").Value = 46000
$ws.Range("I134").Value = 1468.7046
$ws.Range("M134").Value = -1871.1138
$ws.Range("K134").Value = 4406.1138
$ws.Range("L134").Value = 22654.32
$ws.Range("N134").Value = -27724.32
$ws.Range("H134").Value = 3672.5942
$ws.Range("J134").Value = 7551.44

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M16").Value = -3088.111
$ws.Range("K16").Value = 3375.111
$ws.Range("H16").Value = 2850.0588
$ws.Range("I16").Value = 3375.111
$ws.Range("L62").Value = 60000
$ws.Range("H62").Value = 7726.591
$ws.Range("N62").Value = -61248
$ws.Range("J62").Value = 60000
$ws.Range("L65").Value = 300000
$ws.Range("N65").Value = -306240
$ws.Range("H65").Value = 7726.591
$ws.Range("J65").Value = 60000
$ws.Range("M113").Value = -1205.111
$ws.Range("K113").Value = 3375.111
$ws.Range("H113").Value = 2850.0588
$ws.Range("I113").Value = 3375.111
$ws.Range("K122").Value = 125002410
$ws.Range("M122").Value = -124999960
$ws.Range("L122").Value = 25360.5
$ws.Range("N122").Value = -30260.5
$ws.Range("H122").Value = 35716184
$ws.Range("J122").Value = 8453.5
$ws.Range("I122").Value = 41667470

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M80").Value = -1517.8333
$ws.Range("L80").Value = 2209.1667
$ws.Range("K80").Value = 2515.8333
$ws.Range("N80").Value = -4205.1667
$ws.Range("H80").Value = 2439.1667
$ws.Range("J80").Value = 2209.1667
$ws.Range("I80").Value = 2515.8333
$ws.Range("I83").Value = 2515.8333
$ws.Range("M83").Value = -7587.166499999999
$ws.Range("L83").Value = 11045.8335
$ws.Range("H83").Value = 2439.1667
$ws.Range("K83").Value = 12579.1665
$ws.Range("N83").Value = -21029.8335
$ws.Range("J83").Value = 2209.1667
$ws.Range("L122").Value = 20970
$ws.Range("N122").Value = -25870
$ws.Range("H122").Value = 3195
$ws.Range("J122").Value = 6990

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M7").Value = -1769.5
$ws.Range("L7").Value = 1550.8334
$ws.Range("H7").Value = 1791.3182
$ws.Range("K7").Value = 1881.5
$ws.Range("N7").Value = -1774.8334
$ws.Range("J7").Value = 1550.8334
$ws.Range("I7").Value = 1881.5
$ws.Range("M16").Value = -1705.9166
$ws.Range("L16").Value = 50340
$ws.Range("K16").Value = 1875.9166
$ws.Range("N16").Value = -50680
$ws.Range("H16").Value = 8799.357
$ws.Range("J16").Value = 50340
$ws.Range("I16").Value = 1875.9166
$ws.Range("M61").Value = -1285.5
$ws.Range("L61").Value = 9500
$ws.Range("K61").Value = 1487.5
$ws.Range("N61").Value = -9904
$ws.Range("H61").Value = 5493.75
$ws.Range("J61").Value = 9500
$ws.Range("I61").Value = 1487.5
$ws.Range("M113").Value = 682.5
$ws.Range("L113").Value = 9500
$ws.Range("K113").Value = 1487.5
$ws.Range("N113").Value = -13840
$ws.Range("H113").Value = 5493.75
$ws.Range("J113").Value = 9500
$ws.Range("I113").Value = 1487.5
$ws.Range("K126").Value = 5644.5
$ws.Range("M126").Value = -3174.5
$ws.Range("L126").Value = 4652.5002
$ws.Range("N126").Value = -9592.5002
$ws.Range("H126").Value = 1791.3182
$ws.Range("J126").Value = 1550.8334
$ws.Range("I126").Value = 1881.5
$ws.Range("M132").Value = -166847450
$ws.Range("L132").Value = 19424.25
$ws.Range("K132").Value = 166849980
$ws.Range("N132").Value = -24484.25
$ws.Range("H132").Value = 38505836
$ws.Range("J132").Value = 6474.75
$ws.Range("I132").Value = 55616660

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K96").Value = 2869.75
$ws.Range("M96").Value = -1496.75
$ws.Range("L96").Value = 10270
$ws.Range("N96").Value = -13016
$ws.Range("H96").Value = 7579
$ws.Range("J96").Value = 10270
$ws.Range("I96").Value = 2869.75
$ws.Range("K122").Value = 42078.60000000001
$ws.Range("M122").Value = -39628.60000000001
$ws.Range("L122").Value = 12173.5716
$ws.Range("N122").Value = -17073.5716
$ws.Range("H122").Value = 9921.588
$ws.Range("J122").Value = 4057.8572
$ws.Range("I122").Value = 14026.2
$ws.Range("L123").Value = 0
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("N123").ClearContents()
